{"js": "// Replace each two-digit-by-two-digit multiplication prompt in the table\n// with its new pair of operands, keeping the trailing \"=\" and all\n// paragraph/run formatting untouched.\nconst replacements = [\n  [\"21\u00d718=\", \"88\u00d779=\"],\n  [\"27\u00d796=\", \"91\u00d757=\"],\n  [\"53\u00d714=\", \"56\u00d727=\"],\n  [\"46\u00d729=\", \"33\u00d779=\"],\n  [\"73\u00d781=\", \"34\u00d732=\"],\n  [\"30\u00d751=\", \"52\u00d754=\"],\n  [\"40\u00d712=\", \"74\u00d744=\"],\n  [\"84\u00d769=\", \"93\u00d716=\"],\n  [\"37\u00d784=\", \"63\u00d721=\"],\n  [\"41\u00d749=\", \"26\u00d769=\"],\n  [\"94\u00d759=\", \"18\u00d756=\"],\n  [\"70\u00d793=\", \"23\u00d751=\"],\n  [\"47\u00d798=\", \"93\u00d723=\"],\n  [\"31\u00d755=\", \"69\u00d792=\"],\n  [\"26\u00d788=\", \"19\u00d738=\"],\n  [\"43\u00d792=\", \"48\u00d724=\"],\n  [\"34\u00d721=\", \"97\u00d771=\"],\n  [\"39\u00d726=\", \"33\u00d796=\"],\n  [\"60\u00d733=\", \"37\u00d744=\"],\n  [\"29\u00d744=\", \"32\u00d718=\"],\n  [\"59\u00d762=\", \"22\u00d724=\"],\n  [\"34\u00d712=\", \"51\u00d794=\"],\n  [\"39\u00d756=\", \"86\u00d791=\"],\n  [\"28\u00d784=\", \"45\u00d728=\"],\n  [\"95\u00d796=\", \"12\u00d768=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find text to replace: \" + oldText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"21\u00d718=\", \"88\u00d779=\"),\n  @(\"27\u00d796=\", \"91\u00d757=\"),\n  @(\"53\u00d714=\", \"56\u00d727=\"),\n  @(\"46\u00d729=\", \"33\u00d779=\"),\n  @(\"73\u00d781=\", \"34\u00d732=\"),\n  @(\"30\u00d751=\", \"52\u00d754=\"),\n  @(\"40\u00d712=\", \"74\u00d744=\"),\n  @(\"84\u00d769=\", \"93\u00d716=\"),\n  @(\"37\u00d784=\", \"63\u00d721=\"),\n  @(\"41\u00d749=\", \"26\u00d769=\"),\n  @(\"94\u00d759=\", \"18\u00d756=\"),\n  @(\"70\u00d793=\", \"23\u00d751=\"),\n  @(\"47\u00d798=\", \"93\u00d723=\"),\n  @(\"31\u00d755=\", \"69\u00d792=\"),\n  @(\"26\u00d788=\", \"19\u00d738=\"),\n  @(\"43\u00d792=\", \"48\u00d724=\"),\n  @(\"34\u00d721=\", \"97\u00d771=\"),\n  @(\"39\u00d726=\", \"33\u00d796=\"),\n  @(\"60\u00d733=\", \"37\u00d744=\"),\n  @(\"29\u00d744=\", \"32\u00d718=\"),\n  @(\"59\u00d762=\", \"22\u00d724=\"),\n  @(\"34\u00d712=\", \"51\u00d794=\"),\n  @(\"39\u00d756=\", \"86\u00d791=\"),\n  @(\"28\u00d784=\", \"45\u00d728=\"),\n  @(\"95\u00d796=\", \"12\u00d768=\")\n)\n\nforeach ($pair in $replacements) {\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Text = $pair[0]\n  $find.Replacement.ClearFormatting()\n  $find.Replacement.Text = $pair[1]\n  $result = $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n  if (-not $result) {\n    Write-Output \"WARN: replace failed for $($pair[0])\"\n  }\n}\n\n"}
